$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 2227.5715  # ALC!H58: 2085.7334 -> 2227.5715
$ws.Cells.Item(58, 10).Value = 4504.25  # ALC!J58: 3623.4 -> 4504.25
$ws.Cells.Item(58, 12).Value = 13512.75  # ALC!L58: 10870.2 -> 13512.75
$ws.Cells.Item(58, 14).Value = -13812.75  # ALC!N58: -11170.2 -> -13812.75

$ws.Cells.Item(129, 8).Value = 1299.2  # ALC!H129: 1498.6666 -> 1299.2
$ws.Cells.Item(129, 9).Value = 1299.2  # ALC!I129: 1498.6666 -> 1299.2
$ws.Cells.Item(129, 11).Value = 3897.6  # ALC!K129: 4495.9998 -> 3897.6
$ws.Cells.Item(129, 13).Value = 1102.4  # ALC!M129: 504.0002000000004 -> 1102.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1014.25  # ARM!H2: 1005.4 -> 1014.25
$ws.Cells.Item(2, 9).Value = 1014.25  # ARM!I2: 1006 -> 1014.25
$ws.Cells.Item(2, 10).Value = 0  # ARM!J2: 1000 -> 0
$ws.Cells.Item(2, 11).Value = 1014.25  # ARM!K2: 1006 -> 1014.25
$ws.Cells.Item(2, 12).Value = 0  # ARM!L2: 1000 -> 0
$ws.Cells.Item(2, 13).Value = -901.25  # ARM!M2: -893 -> -901.25
$ws.Cells.Item(2, 14).ClearContents()  # ARM!N2: -1226 -> (removed)

$ws.Cells.Item(32, 8).Value = 5114.4707  # ARM!H32: 3369.037 -> 5114.4707
$ws.Cells.Item(32, 9).Value = 5121.625  # ARM!I32: 3306.3076 -> 5121.625
$ws.Cells.Item(32, 11).Value = 5121.625  # ARM!K32: 3306.3076 -> 5121.625
$ws.Cells.Item(32, 13).Value = -4834.625  # ARM!M32: -3019.3076 -> -4834.625

$ws.Cells.Item(45, 8).Value = 2611.1  # ARM!H45: 2410.111 -> 2611.1
$ws.Cells.Item(45, 10).Value = 4437.75  # ARM!J45: 4443.6665 -> 4437.75
$ws.Cells.Item(45, 12).Value = 4437.75  # ARM!L45: 4443.6665 -> 4437.75
$ws.Cells.Item(45, 14).Value = -5191.75  # ARM!N45: -5197.6665 -> -5191.75

$ws.Cells.Item(116, 8).Value = 1014.25  # ARM!H116: 1005.4 -> 1014.25
$ws.Cells.Item(116, 9).Value = 1014.25  # ARM!I116: 1006 -> 1014.25
$ws.Cells.Item(116, 10).Value = 0  # ARM!J116: 1000 -> 0
$ws.Cells.Item(116, 11).Value = 1014.25  # ARM!K116: 1006 -> 1014.25
$ws.Cells.Item(116, 12).Value = 0  # ARM!L116: 1000 -> 0
$ws.Cells.Item(116, 13).Value = 1279.75  # ARM!M116: 1288 -> 1279.75
$ws.Cells.Item(116, 14).ClearContents()  # ARM!N116: -5588 -> (removed)

$ws.Cells.Item(122, 8).Value = 3999  # ARM!H122: 999.5 -> 3999
$ws.Cells.Item(122, 9).Value = 0  # ARM!I122: 999.5 -> 0
$ws.Cells.Item(122, 10).Value = 3999  # ARM!J122: 0 -> 3999
$ws.Cells.Item(122, 11).Value = 0  # ARM!K122: 2998.5 -> 0
$ws.Cells.Item(122, 12).Value = 11997  # ARM!L122: 0 -> 11997
$ws.Cells.Item(122, 13).ClearContents()  # ARM!M122: -548.5 -> (removed)
$ws.Cells.Item(122, 14).Value = -16897  # ARM!N122: None -> -16897

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1014.25  # BSM!H3: 1005.4 -> 1014.25
$ws.Cells.Item(3, 9).Value = 1014.25  # BSM!I3: 1006 -> 1014.25
$ws.Cells.Item(3, 10).Value = 0  # BSM!J3: 1000 -> 0
$ws.Cells.Item(3, 11).Value = 1014.25  # BSM!K3: 1006 -> 1014.25
$ws.Cells.Item(3, 12).Value = 0  # BSM!L3: 1000 -> 0
$ws.Cells.Item(3, 13).Value = -900.25  # BSM!M3: -892 -> -900.25
$ws.Cells.Item(3, 14).ClearContents()  # BSM!N3: -1228 -> (removed)

$ws.Cells.Item(88, 8).Value = 17333  # BSM!H88: 18499.5 -> 17333
$ws.Cells.Item(88, 10).Value = 17333  # BSM!J88: 18499.5 -> 17333
$ws.Cells.Item(88, 12).Value = 17333  # BSM!L88: 18499.5 -> 17333
$ws.Cells.Item(88, 14).Value = -18145  # BSM!N88: -19311.5 -> -18145

$ws.Cells.Item(91, 8).Value = 17333  # BSM!H91: 18499.5 -> 17333
$ws.Cells.Item(91, 10).Value = 17333  # BSM!J91: 18499.5 -> 17333
$ws.Cells.Item(91, 12).Value = 17333  # BSM!L91: 18499.5 -> 17333
$ws.Cells.Item(91, 14).Value = -20141  # BSM!N91: -21307.5 -> -20141

$ws.Cells.Item(94, 8).Value = 1499  # BSM!H94: 0 -> 1499
$ws.Cells.Item(94, 10).Value = 1499  # BSM!J94: 0 -> 1499
$ws.Cells.Item(94, 12).Value = 1499  # BSM!L94: 0 -> 1499
$ws.Cells.Item(94, 14).Value = -2401  # BSM!N94: None -> -2401

$ws.Cells.Item(107, 8).Value = 3121.76  # BSM!H107: 2540.125 -> 3121.76
$ws.Cells.Item(107, 9).Value = 2096.7144  # BSM!I107: 1688.25 -> 2096.7144
$ws.Cells.Item(107, 11).Value = 2096.7144  # BSM!K107: 1688.25 -> 2096.7144
$ws.Cells.Item(107, 13).Value = -176.7143999999998  # BSM!M107: 231.75 -> -176.7143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1749.5  # CRP!H16: 1999 -> 1749.5
$ws.Cells.Item(16, 9).Value = 1749.5  # CRP!I16: 1999 -> 1749.5
$ws.Cells.Item(16, 11).Value = 1749.5  # CRP!K16: 1999 -> 1749.5
$ws.Cells.Item(16, 13).Value = -1462.5  # CRP!M16: -1712 -> -1462.5

$ws.Cells.Item(31, 8).Value = 4382.6665  # CRP!H31: 4024 -> 4382.6665
$ws.Cells.Item(31, 9).Value = 2196.1738  # CRP!I31: 1988.7407 -> 2196.1738
$ws.Cells.Item(31, 11).Value = 2196.1738  # CRP!K31: 1988.7407 -> 2196.1738
$ws.Cells.Item(31, 13).Value = -1901.1738  # CRP!M31: -1693.7407 -> -1901.1738

$ws.Cells.Item(34, 8).Value = 4382.6665  # CRP!H34: 4024 -> 4382.6665
$ws.Cells.Item(34, 9).Value = 2196.1738  # CRP!I34: 1988.7407 -> 2196.1738
$ws.Cells.Item(34, 11).Value = 2196.1738  # CRP!K34: 1988.7407 -> 2196.1738
$ws.Cells.Item(34, 13).Value = -1994.1738  # CRP!M34: -1786.7407 -> -1994.1738

$ws.Cells.Item(58, 8).Value = 3616.0833  # CRP!H58: 3757.6155 -> 3616.0833
$ws.Cells.Item(58, 9).Value = 2575.5  # CRP!I58: 3334.8333 -> 2575.5
$ws.Cells.Item(58, 10).Value = 4656.6665  # CRP!J58: 4120 -> 4656.6665
$ws.Cells.Item(58, 11).Value = 2575.5  # CRP!K58: 3334.8333 -> 2575.5
$ws.Cells.Item(58, 12).Value = 4656.6665  # CRP!L58: 4120 -> 4656.6665
$ws.Cells.Item(58, 13).Value = -2372.5  # CRP!M58: -3131.8333 -> -2372.5
$ws.Cells.Item(58, 14).Value = -5062.6665  # CRP!N58: -4526 -> -5062.6665

$ws.Cells.Item(74, 8).Value = 54254  # CRP!H74: 65382.5 -> 54254
$ws.Cells.Item(74, 10).Value = 54254  # CRP!J74: 65382.5 -> 54254
$ws.Cells.Item(74, 12).Value = 54254  # CRP!L74: 65382.5 -> 54254
$ws.Cells.Item(74, 14).Value = -56002  # CRP!N74: -67130.5 -> -56002

$ws.Cells.Item(77, 8).Value = 54254  # CRP!H77: 65382.5 -> 54254
$ws.Cells.Item(77, 10).Value = 54254  # CRP!J77: 65382.5 -> 54254
$ws.Cells.Item(77, 12).Value = 162762  # CRP!L77: 196147.5 -> 162762
$ws.Cells.Item(77, 14).Value = -171498  # CRP!N77: -204883.5 -> -171498

$ws.Cells.Item(113, 8).Value = 1749.5  # CRP!H113: 1999 -> 1749.5
$ws.Cells.Item(113, 9).Value = 1749.5  # CRP!I113: 1999 -> 1749.5
$ws.Cells.Item(113, 11).Value = 1749.5  # CRP!K113: 1999 -> 1749.5
$ws.Cells.Item(113, 13).Value = 420.5  # CRP!M113: 171 -> 420.5

$ws.Cells.Item(132, 8).Value = 1500  # CRP!H132: 2000 -> 1500
$ws.Cells.Item(132, 9).Value = 1500  # CRP!I132: 2000 -> 1500
$ws.Cells.Item(132, 11).Value = 4500  # CRP!K132: 6000 -> 4500
$ws.Cells.Item(132, 13).Value = -1970  # CRP!M132: -3470 -> -1970

$ws.Cells.Item(134, 8).Value = 2266.0667  # CRP!H134: 2014.3684 -> 2266.0667
$ws.Cells.Item(134, 9).Value = 2266.0667  # CRP!I134: 2014.3684 -> 2266.0667
$ws.Cells.Item(134, 11).Value = 6798.2001  # CRP!K134: 6043.1052 -> 6798.2001
$ws.Cells.Item(134, 13).Value = -4263.2001  # CRP!M134: -3508.1052 -> -4263.2001

$ws.Cells.Item(136, 8).Value = 3616.0833  # CRP!H136: 3757.6155 -> 3616.0833
$ws.Cells.Item(136, 9).Value = 2575.5  # CRP!I136: 3334.8333 -> 2575.5
$ws.Cells.Item(136, 10).Value = 4656.6665  # CRP!J136: 4120 -> 4656.6665
$ws.Cells.Item(136, 11).Value = 7726.5  # CRP!K136: 10004.4999 -> 7726.5
$ws.Cells.Item(136, 12).Value = 13969.9995  # CRP!L136: 12360 -> 13969.9995
$ws.Cells.Item(136, 13).Value = -5176.5  # CRP!M136: -7454.499899999999 -> -5176.5
$ws.Cells.Item(136, 14).Value = -19069.9995  # CRP!N136: -17460 -> -19069.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 198  # CUL!H22: 196.66667 -> 198
$ws.Cells.Item(22, 10).Value = 198  # CUL!J22: 196.66667 -> 198
$ws.Cells.Item(22, 12).Value = 594  # CUL!L22: 590.00001 -> 594
$ws.Cells.Item(22, 14).Value = -932  # CUL!N22: -928.00001 -> -932

$ws.Cells.Item(27, 8).Value = 198  # CUL!H27: 196.66667 -> 198
$ws.Cells.Item(27, 10).Value = 198  # CUL!J27: 196.66667 -> 198
$ws.Cells.Item(27, 12).Value = 594  # CUL!L27: 590.00001 -> 594
$ws.Cells.Item(27, 14).Value = -798  # CUL!N27: -794.00001 -> -798

$ws.Cells.Item(103, 8).Value = 385.375  # CUL!H103: 338.9091 -> 385.375
$ws.Cells.Item(103, 9).Value = 356  # CUL!I103: 301.83334 -> 356
$ws.Cells.Item(103, 10).Value = 414.75  # CUL!J103: 383.4 -> 414.75
$ws.Cells.Item(103, 11).Value = 1068  # CUL!K103: 905.5000200000001 -> 1068
$ws.Cells.Item(103, 12).Value = 1244.25  # CUL!L103: 1150.2 -> 1244.25
$ws.Cells.Item(103, 13).Value = -189  # CUL!M103: -26.50002000000006 -> -189
$ws.Cells.Item(103, 14).Value = -3002.25  # CUL!N103: -2908.2 -> -3002.25

$ws.Cells.Item(113, 8).Value = 1157.8462  # CUL!H113: 1203.1818 -> 1157.8462
$ws.Cells.Item(113, 9).Value = 621  # CUL!I113: 650 -> 621
$ws.Cells.Item(113, 10).Value = 1255.4546  # CUL!J113: 1258.5 -> 1255.4546
$ws.Cells.Item(113, 11).Value = 1863  # CUL!K113: 1950 -> 1863
$ws.Cells.Item(113, 12).Value = 3766.3638  # CUL!L113: 3775.5 -> 3766.3638
$ws.Cells.Item(113, 13).Value = 307  # CUL!M113: 220 -> 307
$ws.Cells.Item(113, 14).Value = -8106.3638  # CUL!N113: -8115.5 -> -8106.3638

$ws.Cells.Item(121, 8).Value = 500  # CUL!H121: 0 -> 500
$ws.Cells.Item(121, 9).Value = 500  # CUL!I121: 0 -> 500
$ws.Cells.Item(121, 11).Value = 1500  # CUL!K121: 0 -> 1500
$ws.Cells.Item(121, 13).Value = -190  # CUL!M121: None -> -190

$ws.Cells.Item(132, 8).Value = 2555.15  # CUL!H132: 3724.5833 -> 2555.15
$ws.Cells.Item(132, 9).Value = 2100.5  # CUL!I132: 2199.5 -> 2100.5
$ws.Cells.Item(132, 10).Value = 2668.8125  # CUL!J132: 4487.125 -> 2668.8125
$ws.Cells.Item(132, 11).Value = 18904.5  # CUL!K132: 19795.5 -> 18904.5
$ws.Cells.Item(132, 12).Value = 24019.3125  # CUL!L132: 40384.125 -> 24019.3125
$ws.Cells.Item(132, 13).Value = -16374.5  # CUL!M132: -17265.5 -> -16374.5
$ws.Cells.Item(132, 14).Value = -29079.3125  # CUL!N132: -45444.125 -> -29079.3125

$ws.Cells.Item(133, 8).Value = 1000  # CUL!H133: 0 -> 1000
$ws.Cells.Item(133, 9).Value = 1000  # CUL!I133: 0 -> 1000
$ws.Cells.Item(133, 11).Value = 3000  # CUL!K133: 0 -> 3000
$ws.Cells.Item(133, 13).Value = 2060  # CUL!M133: None -> 2060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(105, 8).Value = 35000  # GSM!H105: 0 -> 35000
$ws.Cells.Item(105, 10).Value = 35000  # GSM!J105: 0 -> 35000
$ws.Cells.Item(105, 12).Value = 35000  # GSM!L105: 0 -> 35000
$ws.Cells.Item(105, 14).Value = -41988  # GSM!N105: None -> -41988

$ws.Cells.Item(107, 8).Value = 1228.125  # GSM!H107: 1228.1875 -> 1228.125
$ws.Cells.Item(107, 9).Value = 1402.5  # GSM!I107: 1529 -> 1402.5
$ws.Cells.Item(107, 10).Value = 937.5  # GSM!J107: 841.4286 -> 937.5
$ws.Cells.Item(107, 11).Value = 1402.5  # GSM!K107: 1529 -> 1402.5
$ws.Cells.Item(107, 12).Value = 937.5  # GSM!L107: 841.4286 -> 937.5
$ws.Cells.Item(107, 13).Value = 517.5  # GSM!M107: 391 -> 517.5
$ws.Cells.Item(107, 14).Value = -4777.5  # GSM!N107: -4681.4286 -> -4777.5

$ws.Cells.Item(113, 8).Value = 3141.4  # GSM!H113: 5300.1113 -> 3141.4
$ws.Cells.Item(113, 9).Value = 3141.4  # GSM!I113: 2950.3333 -> 3141.4
$ws.Cells.Item(113, 10).Value = 0  # GSM!J113: 9999.666999999999 -> 0
$ws.Cells.Item(113, 11).Value = 3141.4  # GSM!K113: 2950.3333 -> 3141.4
$ws.Cells.Item(113, 12).Value = 0  # GSM!L113: 9999.666999999999 -> 0
$ws.Cells.Item(113, 13).Value = -971.4000000000001  # GSM!M113: -780.3332999999998 -> -971.4000000000001
$ws.Cells.Item(113, 14).ClearContents()  # GSM!N113: -14339.667 -> (removed)

$ws.Cells.Item(132, 8).Value = 3245.9167  # GSM!H132: 3360.4546 -> 3245.9167
$ws.Cells.Item(132, 9).Value = 2696.2  # GSM!I132: 2775.111 -> 2696.2
$ws.Cells.Item(132, 11).Value = 8088.599999999999  # GSM!K132: 8325.332999999999 -> 8088.599999999999
$ws.Cells.Item(132, 13).Value = -5558.599999999999  # GSM!M132: -5795.332999999999 -> -5558.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3395.0667  # LTW!H61: 3763.8462 -> 3395.0667
$ws.Cells.Item(61, 9).Value = 1693.8  # LTW!I61: 1867.75 -> 1693.8
$ws.Cells.Item(61, 11).Value = 1693.8  # LTW!K61: 1867.75 -> 1693.8
$ws.Cells.Item(61, 13).Value = -1491.8  # LTW!M61: -1665.75 -> -1491.8

$ws.Cells.Item(76, 8).Value = 18143.666  # LTW!H76: 18893.75 -> 18143.666
$ws.Cells.Item(76, 10).Value = 18143.666  # LTW!J76: 18893.75 -> 18143.666
$ws.Cells.Item(76, 12).Value = 18143.666  # LTW!L76: 18893.75 -> 18143.666
$ws.Cells.Item(76, 14).Value = -18819.666  # LTW!N76: -19569.75 -> -18819.666

$ws.Cells.Item(79, 8).Value = 18143.666  # LTW!H79: 18893.75 -> 18143.666
$ws.Cells.Item(79, 10).Value = 18143.666  # LTW!J79: 18893.75 -> 18143.666
$ws.Cells.Item(79, 12).Value = 18143.666  # LTW!L79: 18893.75 -> 18143.666
$ws.Cells.Item(79, 14).Value = -20483.666  # LTW!N79: -21233.75 -> -20483.666

$ws.Cells.Item(113, 8).Value = 3395.0667  # LTW!H113: 3763.8462 -> 3395.0667
$ws.Cells.Item(113, 9).Value = 1693.8  # LTW!I113: 1867.75 -> 1693.8
$ws.Cells.Item(113, 11).Value = 1693.8  # LTW!K113: 1867.75 -> 1693.8
$ws.Cells.Item(113, 13).Value = 476.2  # LTW!M113: 302.25 -> 476.2

$ws.Cells.Item(122, 8).Value = 2655.5  # LTW!H122: 2792.3333 -> 2655.5
$ws.Cells.Item(122, 9).Value = 2366.8572  # LTW!I122: 2579.8 -> 2366.8572
$ws.Cells.Item(122, 11).Value = 7100.571599999999  # LTW!K122: 7739.400000000001 -> 7100.571599999999
$ws.Cells.Item(122, 13).Value = -4650.571599999999  # LTW!M122: -5289.400000000001 -> -4650.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 0  # WVR!H69: 2500 -> 0
$ws.Cells.Item(69, 10).Value = 0  # WVR!J69: 2500 -> 0
$ws.Cells.Item(69, 12).Value = 0  # WVR!L69: 2500 -> 0
$ws.Cells.Item(69, 14).ClearContents()  # WVR!N69: -3998 -> (removed)

$ws.Cells.Item(72, 8).Value = 0  # WVR!H72: 2500 -> 0
$ws.Cells.Item(72, 10).Value = 0  # WVR!J72: 2500 -> 0
$ws.Cells.Item(72, 12).Value = 0  # WVR!L72: 7500 -> 0
$ws.Cells.Item(72, 14).ClearContents()  # WVR!N72: -14988 -> (removed)

$ws.Cells.Item(82, 8).Value = 12000  # WVR!H82: 0 -> 12000
$ws.Cells.Item(82, 10).Value = 12000  # WVR!J82: 0 -> 12000
$ws.Cells.Item(82, 12).Value = 12000  # WVR!L82: 0 -> 12000
$ws.Cells.Item(82, 14).Value = -12766  # WVR!N82: None -> -12766

$ws.Cells.Item(85, 8).Value = 12000  # WVR!H85: 0 -> 12000
$ws.Cells.Item(85, 10).Value = 12000  # WVR!J85: 0 -> 12000
$ws.Cells.Item(85, 12).Value = 12000  # WVR!L85: 0 -> 12000
$ws.Cells.Item(85, 14).Value = -14652  # WVR!N85: None -> -14652

$ws.Cells.Item(107, 8).Value = 1399.1  # WVR!H107: 1405.5 -> 1399.1
$ws.Cells.Item(107, 9).Value = 1734.4286  # WVR!I107: 1967.6666 -> 1734.4286
$ws.Cells.Item(107, 10).Value = 616.6667  # WVR!J107: 562.25 -> 616.6667
$ws.Cells.Item(107, 11).Value = 5203.2858  # WVR!K107: 5902.9998 -> 5203.2858
$ws.Cells.Item(107, 12).Value = 1850.0001  # WVR!L107: 1686.75 -> 1850.0001
$ws.Cells.Item(107, 13).Value = -3283.2858  # WVR!M107: -3982.9998 -> -3283.2858
$ws.Cells.Item(107, 14).Value = -5690.0001  # WVR!N107: -5526.75 -> -5690.0001

$ws.Cells.Item(132, 8).Value = 2749.5  # WVR!H132: 2999 -> 2749.5
$ws.Cells.Item(132, 9).Value = 2500  # WVR!I132: 0 -> 2500
$ws.Cells.Item(132, 11).Value = 7500  # WVR!K132: 0 -> 7500
$ws.Cells.Item(132, 13).Value = -4970  # WVR!M132: None -> -4970
